# Update "想去人数" (interested-count) values in column F for the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets, as produced by the
# gh-pages data regeneration commit.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 371
    4  = 10828
    8  = 1339
    9  = 8297
    15 = 3304
    18 = 31
    19 = 787
    20 = 132
    24 = 1781
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
